# Update Data Sources from LFX
# Re-style every table that currently uses the old "Data Sources from LFX"
# table style so it instead uses the new table style GUID.

$oldStyleId = "{50C352CE-1B1D-4D2E-B4E8-68B4B13091BA}"
$newStyleId = "{27CCD006-CB68-44EF-8F73-6356A743C068}"

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)

        $hasTable = $false
        try { $hasTable = [bool]$sh.HasTable } catch { $hasTable = $false }

        if ($hasTable) {
            $tbl = $sh.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
